{"js": "// Replace the 15 lattice-multiplication problems in the single 5x3 table\n// with a new set of problems, keeping each cell's structure (one run,\n// sz=32, line breaks between the 5 lines) intact.\n//\n// Each cell holds 5 lines separated by manual line breaks (<w:br/>, i.e. a\n// vertical-tab \"\\v\" in the Office.js text model):\n//   1) \"AB x CD\"        -- the problem\n//   2) \"  C    D\"       -- spaced digits of the second factor\n//   3) \"  ----\"         -- separator\n//   4) \"A|    |\"        -- first digit of the first factor\n//   5) \"B|    |\"        -- second digit of the first factor\n\nconst newCells = [\n  [\"38 x 15\", \"  1    5\", \"  ----\", \"3|    |\", \"8|    |\"],\n  [\"16 x 92\", \"  9    2\", \"  ----\", \"1|    |\", \"6|    |\"],\n  [\"12 x 99\", \"  9    9\", \"  ----\", \"1|    |\", \"2|    |\"],\n  [\"84 x 69\", \"  6    9\", \"  ----\", \"8|    |\", \"4|    |\"],\n  [\"56 x 16\", \"  1    6\", \"  ----\", \"5|    |\", \"6|    |\"],\n  [\"46 x 74\", \"  7    4\", \"  ----\", \"4|    |\", \"6|    |\"],\n  [\"73 x 30\", \"  3    0\", \"  ----\", \"7|    |\", \"3|    |\"],\n  [\"54 x 10\", \"  1    0\", \"  ----\", \"5|    |\", \"4|    |\"],\n  [\"50 x 17\", \"  1    7\", \"  ----\", \"5|    |\", \"0|    |\"],\n  [\"71 x 49\", \"  4    9\", \"  ----\", \"7|    |\", \"1|    |\"],\n  [\"97 x 24\", \"  2    4\", \"  ----\", \"9|    |\", \"7|    |\"],\n  [\"65 x 30\", \"  3    0\", \"  ----\", \"6|    |\", \"5|    |\"],\n  [\"78 x 55\", \"  5    5\", \"  ----\", \"7|    |\", \"8|    |\"],\n  [\"23 x 86\", \"  8    6\", \"  ----\", \"2|    |\", \"3|    |\"],\n  [\"20 x 54\", \"  5    4\", \"  ----\", \"2|    |\", \"0|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst colCount = table.values[0].length;\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    const lines = newCells[idx];\n    idx++;\n    para.insertText(lines.join(\"\\v\"), Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication problems in the single 5x3 table\n# with a new set of problems, keeping each cell's structure (one run,\n# sz=32, manual line breaks between the 5 lines) intact.\n#\n# Each cell holds 5 lines separated by manual line breaks (a vertical-tab,\n# \"`v\", maps to <w:br/> when written through Range.Text):\n#   1) \"AB x CD\"        -- the problem\n#   2) \"  C    D\"       -- spaced digits of the second factor\n#   3) \"  ----\"         -- separator\n#   4) \"A|    |\"        -- first digit of the first factor\n#   5) \"B|    |\"        -- second digit of the first factor\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newCells = @(\n    @(\"38 x 15\", \"  1    5\", \"  ----\", \"3|    |\", \"8|    |\"),\n    @(\"16 x 92\", \"  9    2\", \"  ----\", \"1|    |\", \"6|    |\"),\n    @(\"12 x 99\", \"  9    9\", \"  ----\", \"1|    |\", \"2|    |\"),\n    @(\"84 x 69\", \"  6    9\", \"  ----\", \"8|    |\", \"4|    |\"),\n    @(\"56 x 16\", \"  1    6\", \"  ----\", \"5|    |\", \"6|    |\"),\n    @(\"46 x 74\", \"  7    4\", \"  ----\", \"4|    |\", \"6|    |\"),\n    @(\"73 x 30\", \"  3    0\", \"  ----\", \"7|    |\", \"3|    |\"),\n    @(\"54 x 10\", \"  1    0\", \"  ----\", \"5|    |\", \"4|    |\"),\n    @(\"50 x 17\", \"  1    7\", \"  ----\", \"5|    |\", \"0|    |\"),\n    @(\"71 x 49\", \"  4    9\", \"  ----\", \"7|    |\", \"1|    |\"),\n    @(\"97 x 24\", \"  2    4\", \"  ----\", \"9|    |\", \"7|    |\"),\n    @(\"65 x 30\", \"  3    0\", \"  ----\", \"6|    |\", \"5|    |\"),\n    @(\"78 x 55\", \"  5    5\", \"  ----\", \"7|    |\", \"8|    |\"),\n    @(\"23 x 86\", \"  8    6\", \"  ----\", \"2|    |\", \"3|    |\"),\n    @(\"20 x 54\", \"  5    4\", \"  ----\", \"2|    |\", \"0|    |\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $lines = $newCells[$idx]\n        $idx++\n        $cell.Range.Text = ($lines -join \"`v\")\n    }\n}\n"}
